$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'workout knee pad'
$ws.Range('A2').Value = 'girls knee pads volleyball'
$ws.Range('A3').Value = 'youth indoor volleyball'
$ws.Range('A4').Value = 'boys compressions'
$ws.Range('A5').Value = 'knee sleeve for bursitis'
$ws.Range('A6').Value = 'stretch marks men'
$ws.Range('A7').Value = 'mens below the knee shorts'
$ws.Range('A8').Value = 'sport tights men'
$ws.Range('A9').Value = 'joint compression'
$ws.Range('A10').Value = 'softball sweats'
$ws.Range('A11').Value = 'youth sport leggings for boys'
$ws.Range('A12').Value = 'padded compression knee sleeve'
$ws.Range('A13').Value = 'youth volleyball shorts'
$ws.Range('A14').Value = 'knee hockey'
$ws.Range('A15').Value = 'girls hockey compression shorts'
$ws.Range('A16').Value = 'football pads adult'
$ws.Range('A17').Value = 'padded football shorts'
$ws.Range('A18').Value = 'basketball shorts for men xxl'
$ws.Range('A19').Value = 'womens softball compression shorts'
$ws.Range('A20').Value = 'basketball leg sleeves youth boys'
$ws.Range('A21').Value = 'football youth tights'
$ws.Range('A22').Value = 'paintball protective gear pants'
$ws.Range('A23').Value = 'catcher pads'
$ws.Range('A24').Value = 'size chart for men'
$ws.Range('A25').Value = 'paintball shorts padded'
$ws.Range('A26').Value = 'men spandex leggings'
$ws.Range('A27').Value = 'calf compression sleeve for youth'
$ws.Range('A28').Value = 'mens 3/4 shorts'
$ws.Range('A29').Value = 'tendinitis knee'
$ws.Range('A30').Value = 'basketball hopes'
$ws.Range('A31').Value = 'mens volleyball shorts'
$ws.Range('A32').Value = 'cycling pants men padded'
$ws.Range('A33').Value = 'knee guards volleyball'
$ws.Range('A34').Value = 'silicone calf pads'
$ws.Range('A35').Value = 'work knees pads'
$ws.Range('A36').Value = 'youth padded knee sleeves'
$ws.Range('A37').Value = 'knee protector running'
$ws.Range('A38').Value = 'sliding shorts youth boys'
$ws.Range('A39').Value = 'padded leg sleeve'
$ws.Range('A40').Value = 'snowboarding pants for men'
$ws.Range('A41').Value = 'basketball sweat pants for men'
$ws.Range('A42').Value = 'black baseball pants youth boys'
$ws.Range('A43').Value = '1/2x28 thread protector'
$ws.Range('A44').Value = 'six one knee pads'
$ws.Range('A45').Value = 'outdoor knee pads'
$ws.Range('A46').Value = 'shorts with leggings men'
$ws.Range('A47').Value = 'youth baseball compression'
$ws.Range('A48').Value = 'hip guards for fall protection'
$ws.Range('A49').Value = 'mens compression yoga pants'
$ws.Range('A50').Value = 'basketball under shorts'
$ws.Range('A51').Value = 'capris spandex'
$ws.Range('A52').Value = 'men leg compression pants'
$ws.Range('A53').Value = 'youth compression padded shorts'
$ws.Range('A54').Value = 'youth volleyball knee pads for girls'
$ws.Range('A55').Value = 'baseball short pants'
$ws.Range('A56').Value = 'youth sports leggings boys'
$ws.Range('A57').Value = 'youth athletic tights boys'
$ws.Range('A58').Value = 'compression pant for men'
$ws.Range('A59').Value = 'working knee pad'
$ws.Range('A60').Value = 'mens long shorts below knee'
$ws.Range('A61').Value = 'black youth football pants'
$ws.Range('A62').Value = 'baseball pants youth xxl'
$ws.Range('A63').Value = 'girls compression leggings'
$ws.Range('A64').Value = 'easy knee pads'
$ws.Range('A65').Value = 'compression knee sleeve youth'
$ws.Range('A66').Value = 'basketball gear for boys youth'
$ws.Range('A67').Value = 'men leggings black'
$ws.Range('A68').Value = 'boy leggings for sports youth'
$ws.Range('A69').Value = 'leggings men compression'
$ws.Range('A70').Value = 'hockey pads youth'
$ws.Range('A71').Value = 'under knee pad sleeves'
$ws.Range('A72').Value = 'weightlifting floor'
$ws.Range('A73').Value = 'anti sweat pads'
$ws.Range('A74').Value = 'mens knee shorts'
$ws.Range('A75').Value = 'leggings for sports men'
$ws.Range('A76').Value = 'girls knee pads volleyball youth'
$ws.Range('A77').Value = 'compression knee sleeve with pad'
$ws.Range('A78').Value = 'running pad'
$ws.Range('A79').Value = 'baseball stretch bands'
$ws.Range('A80').Value = 'softball long pants'
$ws.Range('A81').Value = 'elastic waisted pants for men'
$ws.Range('A82').Value = 'mountain bike knee pads for men'
$ws.Range('A83').Value = 'yoga pant for men'
$ws.Range('A84').Value = 'compression basketball knee sleeve'
$ws.Range('A85').Value = 'compression knee sleeves with pads'
$ws.Range('A86').Value = 'knee work pad'
$ws.Range('A87').Value = 'running compression pants'
$ws.Range('A88').Value = 'professional knee pads construction'
$ws.Range('A89').Value = 'compression pants for girls'
$ws.Range('A90').Value = 'womans softball sliding shorts'
$ws.Range('A91').Value = 'compression calf sleeve youth'
$ws.Range('A92').Value = 'performance compression knee sleeve'
$ws.Range('A93').Value = 'youth wrestling kneepads'
$ws.Range('A94').Value = 'knee pads for work women'
$ws.Range('A95').Value = 'baseball pants short'
$ws.Range('A96').Value = 'men capri shorts'
$ws.Range('A97').Value = 'volleyball knee pads girls'
$ws.Range('A98').Value = 'sit pad hiking'
$ws.Range('A99').Value = 'tactical pants knee pads'
$ws.Range('A100').Value = 'training pants men'
